$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 10448
$ws.Range("E2").Value = 1152
$ws.Range("F2").Value = 1152
$ws.Range("G2").Value = 1401
$ws.Range("H2").Value = 1015
$ws.Range("I2").Value = 568
$ws.Range("J2").Value = 447
$ws.Range("K2").Value = 17909
$ws.Range("L2").Value = 6098
$ws.Range("M2").Value = 11811
$ws.Range("N2").Value = 7047
$ws.Range("O2").Value = 4764
$ws.Range("P2").Value = 266
$ws.Range("Q2").Value = 545
$ws.Range("R2").Value = -1029
$ws.Range("S2").Value = 553
$ws.Range("T2").Value = 783
$ws.Range("U2").Value = -237
$ws.Range("V2").Value = 2550
$ws.Range("W2").Value = 11.03
$ws.Range("X2").Value = 9.710000000000001
$ws.Range("Y2").Value = 8.24
$ws.Range("Z2").Value = 5.94
$ws.Range("AA2").Value = 51.63
$ws.Range("AB2").Value = 2455.16
$ws.Range("AC2").Value = 1146
$ws.Range("AD2").Value = 18.24
$ws.Range("AE2").Value = 15455
$ws.Range("AF2").Value = 1.35
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 1.44
$ws.Range("AI2").Value = 24.1
$ws.Range("AJ2").Value = 47028210

# Row 3
$ws.Range("D3").Value = 11329
$ws.Range("E3").Value = 1076
$ws.Range("F3").Value = 1076
$ws.Range("G3").Value = 1570
$ws.Range("H3").Value = 1128
$ws.Range("I3").Value = 654
$ws.Range("J3").Value = 474
$ws.Range("K3").Value = 20587
$ws.Range("L3").Value = 7513
$ws.Range("M3").Value = 13074
$ws.Range("N3").Value = 7322
$ws.Range("O3").Value = 5752
$ws.Range("P3").Value = 266
$ws.Range("Q3").Value = 70
$ws.Range("R3").Value = -201
$ws.Range("S3").Value = 1263
$ws.Range("T3").Value = 1518
$ws.Range("U3").Value = -1448
$ws.Range("V3").Value = 3380
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 9.960000000000001
$ws.Range("Y3").Value = 9.109999999999999
$ws.Range("Z3").Value = 5.86
$ws.Range("AA3").Value = 57.46
$ws.Range("AB3").Value = 2649.77
$ws.Range("AC3").Value = 1320
$ws.Range("AD3").Value = 31.77
$ws.Range("AE3").Value = 16127
$ws.Range("AF3").Value = 2.6
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 0.72
$ws.Range("AI3").Value = 20.82
$ws.Range("AJ3").Value = 47028210

# Row 4
$ws.Range("D4").Value = 13545
$ws.Range("E4").Value = 993
$ws.Range("F4").Value = 993
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 743
$ws.Range("I4").Value = 472
$ws.Range("J4").Value = 271
$ws.Range("K4").Value = 22493
$ws.Range("L4").Value = 8723
$ws.Range("M4").Value = 13770
$ws.Range("N4").Value = 7716
$ws.Range("O4").Value = 6054
$ws.Range("P4").Value = 266
$ws.Range("Q4").Value = -8
$ws.Range("R4").Value = -1995
$ws.Range("S4").Value = 1438
$ws.Range("T4").Value = 1599
$ws.Range("U4").Value = -1607
$ws.Range("V4").Value = 4910
$ws.Range("W4").Value = 7.33
$ws.Range("X4").Value = 5.48
$ws.Range("Y4").Value = 6.27
$ws.Range("Z4").Value = 3.45
$ws.Range("AA4").Value = 63.34
$ws.Range("AB4").Value = 2792.61
$ws.Range("AC4").Value = 952
$ws.Range("AD4").Value = 24.99
$ws.Range("AE4").Value = 16995
$ws.Range("AF4").Value = 1.4
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 1.26
$ws.Range("AI4").Value = 28.87
$ws.Range("AJ4").Value = 47028210

# Row 5
$ws.Range("D5").Value = 14655
$ws.Range("E5").Value = 1120
$ws.Range("F5").Value = 1120
$ws.Range("G5").Value = 961
$ws.Range("H5").Value = 707
$ws.Range("I5").Value = 434
$ws.Range("J5").Value = 274
$ws.Range("K5").Value = 25607
$ws.Range("L5").Value = 11210
$ws.Range("M5").Value = 14397
$ws.Range("N5").Value = 8128
$ws.Range("O5").Value = 6269
$ws.Range("P5").Value = 266
$ws.Range("Q5").Value = 1445
$ws.Range("R5").Value = -1573
$ws.Range("S5").Value = 1835
$ws.Range("T5").Value = 1749
$ws.Range("U5").Value = -304
$ws.Range("V5").Value = 6916
$ws.Range("W5").Value = 7.64
$ws.Range("X5").Value = 4.83
$ws.Range("Y5").Value = 5.47
$ws.Range("Z5").Value = 2.94
$ws.Range("AA5").Value = 77.87
$ws.Range("AB5").Value = 2901.96
$ws.Range("AC5").Value = 875
$ws.Range("AD5").Value = 44.84
$ws.Range("AE5").Value = 17903
$ws.Range("AF5").Value = 2.19
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 0.76
$ws.Range("AI5").Value = 31.41
$ws.Range("AJ5").Value = 47028210

# Row 6
$ws.Range("D6").Value = 15487
$ws.Range("E6").Value = 493
$ws.Range("F6").Value = 493
$ws.Range("G6").Value = 432
$ws.Range("H6").Value = 194
$ws.Range("I6").Value = 185
$ws.Range("K6").Value = 26286
$ws.Range("L6").Value = 11755
$ws.Range("M6").Value = 14531
$ws.Range("N6").Value = 8043
$ws.Range("P6").Value = 266
$ws.Range("Q6").Value = 376
$ws.Range("R6").Value = -1717
$ws.Range("S6").Value = 1146
$ws.Range("T6").Value = 1180
$ws.Range("U6").Value = -804
$ws.Range("V6").Value = 7925
$ws.Range("W6").Value = 3.19
$ws.Range("X6").Value = 1.25
$ws.Range("Y6").Value = 2.28
$ws.Range("Z6").Value = 0.75
$ws.Range("AA6").Value = 80.89
$ws.Range("AB6").Value = 3004.07
$ws.Range("AC6").Value = 373
$ws.Range("AD6").Value = 67.34
$ws.Range("AE6").Value = 17716
$ws.Range("AH6").Value = 1.42
$ws.Range("AI6").Value = 61.47
$ws.Range("AJ6").Value = 47028210

# Clear AF6:AG6 entirely (cells removed)
$ws.Range("AF6:AG6").ClearContents()

# Clear data cells for rows 7-9 (only A/B/C remain)
$ws.Range("D7:AJ9").ClearContents()
